$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.392.93"
$ws.Range("E2").Value = "  -1.25%  "

$ws.Range("D3").Value = "2.519.11"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'316.54"
$ws.Range("E5").Value = "  +3.57%  "

$ws.Range("D6").Value = "'94.05"
$ws.Range("E6").Value = "  -7.91%  "

$ws.Range("E7").Value = "  -0.79%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.525"
$ws.Range("E9").Value = "  -3.84%  "

$ws.Range("D10").Value = "'35.66"
$ws.Range("E10").Value = "  -6.04%  "

$ws.Range("D11").Value = "'0.0805"
$ws.Range("E11").Value = "  -1.29%  "

$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").Value = "'7.49"
$ws.Range("E13").Value = "  -3.31%  "

$ws.Range("D14").Value = "2.905.76"
$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.537.98"
$ws.Range("E15").Value = "  +1.34%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'15.31"
$ws.Range("E16").Value = "  +0.82%  "

$ws.Range("E17").Value = "  -3.82%  "

$ws.Range("D18").Value = "42.460.07"
$ws.Range("E18").Value = "  -1.13%  "

$ws.Range("D19").Value = "'13.01"
$ws.Range("E19").Value = "  -1.39%  "

$ws.Range("D20").Value = "'6.53"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("D21").Value = "0.0₃0956"
$ws.Range("E21").Value = "  -3.05%  "

$ws.Range("D22").Value = "'69.97"
$ws.Range("E22").Value = "  -2.32%  "

$ws.Range("D23").Value = "'250.44"
$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("D24").Value = "'2.94"
$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("E25").Value = "  -2.93%  "

$ws.Range("D26").Value = "'26.33"
$ws.Range("E26").Value = "  -3.21%  "

$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("E28").Value = "  +2.18%  "

$ws.Range("E29").Value = "  -1.99%  "

$ws.Range("D30").Value = "'38.73"
$ws.Range("E30").Value = "  -1.25%  "

$ws.Range("D31").Value = "'5.90"
$ws.Range("E31").Value = "  -4.53%  "

$ws.Range("D32").Value = "'155.43"
$ws.Range("E32").Value = "  -1.43%  "

$ws.Range("D33").Value = "'19.03"
$ws.Range("E33").Value = "  +3.36%  "

$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").Value = "'0.0781"
$ws.Range("E36").Value = "  -2.00%  "

$ws.Range("D37").Value = "'2.63"
$ws.Range("E37").Value = "  -0.96%  "

$ws.Range("E38").Value = "  -4.46%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.118"
$ws.Range("E39").Value = "  -1.05%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'23.54"
$ws.Range("E40").Value = "  -2.56%  "

$ws.Range("E41").Value = "  +9.99%  "

$ws.Range("E42").Value = "  +0.40%  "

$ws.Range("D43").Value = "'3.76"
$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("E44").Value = "  -2.21%  "

$ws.Range("D45").Value = "'3.26"
$ws.Range("E45").Value = "  -6.80%  "

$ws.Range("D46").Value = "2.003.31"
$ws.Range("E46").Value = "  -2.09%  "

$ws.Range("D47").Value = "'84.57"
$ws.Range("E47").Value = "  -2.16%  "

$ws.Range("D48").Value = "'8.79"
$ws.Range("E48").Value = "  -1.93%  "

$ws.Range("D49").Value = "2.761.40"
$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("D50").Value = "'73.70"
$ws.Range("E50").Value = "  +1.11%  "

$ws.Range("D51").Value = "'101.62"
$ws.Range("E51").Value = "  -1.34%  "

